# Plan to use MBO.xlsx - Add files via upload
# Rewrites the item list (rows 4+) with the new set of budget line items,
# adds new shared strings, and updates dimension/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Clear out the old data rows (4 through old last row 31) in columns
#    D:H so we can rewrite the whole item table cleanly.
# ---------------------------------------------------------------------
$ws.Range("D4:H31").Clear()

# ---------------------------------------------------------------------
# 2. Helper data: each row is Item, USD(qty), VND(fixed price), HasBuy(x)
#    HasBuy = $true writes "x" into column D for that row.
# ---------------------------------------------------------------------
$rows = @(
    @{ Item = "Glass";                                               G = 600000;  Buy = $true  },
    @{ Item = "Sun Glass";                                           G = 2000000; Buy = $false },
    @{ Item = "Nhổ răng";                                            G = 2000000; Buy = $true  },
    @{ Item = "Chair + monitor stand";                               G = 1700000; Buy = $true  },
    @{ Item = "Cherry G80";                                          G = 2650000; Buy = $true  },
    @{ Item = "Cover Cherry G80";                                    G = 450000;  Buy = $true  },
    @{ Item = "ELECOM M-HT1DRBK Wireless Trackball Mouse";           G = 2150000; Buy = $false },
    @{ Item = "Cover for Elecom Trackball Mouse";                    G = 700000;  Buy = $false },
    @{ Item = "ELECOM M-HT1URBK Wired Trackball Mouse";              G = 1700000; Buy = $false },
    @{ Item = "Watch band";                                          G = 1000000; Buy = $true  },
    @{ Item = "Vacation";                                            G = 3000000; Buy = $false },
    @{ Item = "Wrist Rest (from Vietnam)";                           G = 400000;  Buy = $true  },
    @{ Item = "Memory card";                                         G = 600000;  Buy = $false },
    @{ Item = "HDD External";                                        G = 4100000; Buy = $false },
    @{ Item = "CD/DVD External";                                     G = 750000;  Buy = $false },
    @{ Item = "Balô";                                                G = 500000;  Buy = $false },
    @{ Item = "Whey";                                                G = 1350000; Buy = $true  },
    @{ Item = "Vitamin";                                             G = 850000;  Buy = $true  },
    @{ Item = "BCAA";                                                G = 750000;  Buy = $true  },
    @{ Item = "SMD Gateron";                                         F = 22;      Buy = $false },
    @{ Item = "Enjoyt PBT Blank color1";                             F = 38.9;    Buy = $false },
    @{ Item = "Enjoyt PBT Korean & Blue";                            F = 89;      Buy = $false },
    @{ Item = "Enjoyt PBT Russian & Green";                          F = 89;      Buy = $false },
    @{ Item = "Enjoy PBT Japanese";                                  F = 89;      Buy = $false },
    @{ Item = "BLACK DYE-SUBLIMATED PBT KEYSET (Korean 151)";        F = 63;      Buy = $false },
    @{ Item = "Front/Side Printed Backlit Keycaps (ANSI 104)";       F = 35;      Buy = $false },
    @{ Item = "WASD Code White";                                     F = 150;     Buy = $false },
    @{ Item = "WASD Code Black 2nd Hand";                            F = 100;     Buy = $false },
    @{ Item = "WASD White (No keycap)";                              F = 100;     Buy = $false },
    @{ Item = "Keyreative - Vulcan PBT Keyset (Standard)";           F = 65;      Buy = $false },
    @{ Item = 'DSA "Think Different" Keyset (Base)';                 F = 70;      Buy = $false },
    @{ Item = "GMK Corsa Auto - Base+Novelties";                     F = 199;     Buy = $false },
    @{ Item = "GMK Grìseann";                                        F = 175;     Buy = $false },
    @{ Item = "JTK White on Black (Full, Base + Mod)";               F = 79;      Buy = $false },
    @{ Item = "JTK Red on White (Full, Base + Mod)";                 F = 105;     Buy = $false },
    @{ Item = "JTK Red on Black (Base)";                             F = 72;      Buy = $false },
    @{ Item = "Keyreative Gradient PBT Keyset - Ocean";              F = 52;      Buy = $false },
    @{ Item = "Keyreative - Klingon PBT Keyset (Standard)";          F = 85;      Buy = $false },
    @{ Item = "Keyreative - Vulcan PBT Keyset (Glow in the Dark)";   F = 94;      Buy = $false },
    @{ Item = "Keyreative - Klingon PBT Keyset (Glow in the Dark)";  F = 78;      Buy = $false },
    @{ Item = "Keyreative - Klingon PBT Keyset (Standard)";          F = 88;      Buy = $false }
)

$r = 4
foreach ($row in $rows) {
    if ($row.Buy) {
        $ws.Cells.Item($r, 4).Value = "x"
    }
    $ws.Cells.Item($r, 5).Value = $row.Item
    if ($row.ContainsKey("F")) {
        $ws.Cells.Item($r, 6).Value = $row.F
    } else {
        $ws.Cells.Item($r, 7).Value = $row.G
    }
    $ws.Cells.Item($r, 8).Formula = "=IF(F$r=0,G$r,F$r*25000)"
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Re-apply the sheet view (scroll position + active cell) and let the
#    dimension follow the new used range automatically.
# ---------------------------------------------------------------------
$ws.Range("E41").Select() | Out-Null

$wb.Save()
